{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// expression in the worksheet table with the updated values from the\n// target revision. Every \"old\" string in the table (including the date\n// paragraph) is unique within the document, so a simple\n// search-and-replace per pair is unambiguous and non-overlapping.\nconst replacements = [\n  [\"2024-07-19 Friday\", \"2024-07-20 Saturday\"],\n  [\"73\u00d740=\", \"34\u00d740=\"],\n  [\"70\u00d736=\", \"29\u00d776=\"],\n  [\"77\u00d723=\", \"16\u00d729=\"],\n  [\"99\u00d786=\", \"31\u00d722=\"],\n  [\"52\u00d736=\", \"53\u00d751=\"],\n  [\"71\u00d728=\", \"25\u00d739=\"],\n  [\"40\u00d727=\", \"21\u00d774=\"],\n  [\"53\u00d720=\", \"86\u00d764=\"],\n  [\"78\u00d746=\", \"42\u00d792=\"],\n  [\"29\u00d729=\", \"11\u00d748=\"],\n  [\"24\u00d785=\", \"63\u00d786=\"],\n  [\"65\u00d755=\", \"97\u00d780=\"],\n  [\"76\u00d764=\", \"38\u00d780=\"],\n  [\"47\u00d720=\", \"41\u00d730=\"],\n  [\"13\u00d757=\", \"16\u00d794=\"],\n  [\"33\u00d722=\", \"15\u00d778=\"],\n  [\"85\u00d799=\", \"31\u00d781=\"],\n  [\"62\u00d758=\", \"81\u00d769=\"],\n  [\"91\u00d769=\", \"37\u00d791=\"],\n  [\"34\u00d743=\", \"25\u00d757=\"],\n  [\"87\u00d776=\", \"61\u00d751=\"],\n  [\"88\u00d748=\", \"68\u00d778=\"],\n  [\"21\u00d785=\", \"77\u00d799=\"],\n  [\"96\u00d762=\", \"15\u00d763=\"],\n  [\"37\u00d747=\", \"58\u00d730=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and every two-digit x two-digit multiplication\n# prompt in the table to match the target revision. Every \"old\" string is\n# unique within the document (verified against the source OOXML), so a\n# plain Find/Replace (MatchCase, no wildcards, whole-story scope) per pair\n# is unambiguous.\n$pairs = @(\n  @(\"2024-07-19 Friday\", \"2024-07-20 Saturday\"),\n  @(\"73\u00d740=\", \"34\u00d740=\"),\n  @(\"70\u00d736=\", \"29\u00d776=\"),\n  @(\"77\u00d723=\", \"16\u00d729=\"),\n  @(\"99\u00d786=\", \"31\u00d722=\"),\n  @(\"52\u00d736=\", \"53\u00d751=\"),\n  @(\"71\u00d728=\", \"25\u00d739=\"),\n  @(\"40\u00d727=\", \"21\u00d774=\"),\n  @(\"53\u00d720=\", \"86\u00d764=\"),\n  @(\"78\u00d746=\", \"42\u00d792=\"),\n  @(\"29\u00d729=\", \"11\u00d748=\"),\n  @(\"24\u00d785=\", \"63\u00d786=\"),\n  @(\"65\u00d755=\", \"97\u00d780=\"),\n  @(\"76\u00d764=\", \"38\u00d780=\"),\n  @(\"47\u00d720=\", \"41\u00d730=\"),\n  @(\"13\u00d757=\", \"16\u00d794=\"),\n  @(\"33\u00d722=\", \"15\u00d778=\"),\n  @(\"85\u00d799=\", \"31\u00d781=\"),\n  @(\"62\u00d758=\", \"81\u00d769=\"),\n  @(\"91\u00d769=\", \"37\u00d791=\"),\n  @(\"34\u00d743=\", \"25\u00d757=\"),\n  @(\"87\u00d776=\", \"61\u00d751=\"),\n  @(\"88\u00d748=\", \"68\u00d778=\"),\n  @(\"21\u00d785=\", \"77\u00d799=\"),\n  @(\"96\u00d762=\", \"15\u00d763=\"),\n  @(\"37\u00d747=\", \"58\u00d730=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Could not find expected text: $oldText\"\n  }\n}\n"}
